$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (item/gene changed)
$ws.Cells.Item(2, 1).Value = "cg01620164"
$ws.Cells.Item(2, 2).Value = "FIGN"
$ws.Cells.Item(2, 3).Value = "auto"
$ws.Cells.Item(2, 4).Value = "auto"

# Row 3
$ws.Cells.Item(3, 1).Value = "cg14079463"
$ws.Cells.Item(3, 2).Value = "C6orf174"
$ws.Cells.Item(3, 3).Value = "auto"
$ws.Cells.Item(3, 4).Value = "auto"

# Row 5
$ws.Cells.Item(5, 1).Value = "cg17076667"
$ws.Cells.Item(5, 2).Value = "BAG3"
$ws.Cells.Item(5, 3).Value = "auto"
$ws.Cells.Item(5, 4).Value = "auto"

# Row 6
$ws.Cells.Item(6, 1).Value = "cg27615582"
$ws.Cells.Item(6, 2).Value = "PRR4"
$ws.Cells.Item(6, 3).Value = "auto"
$ws.Cells.Item(6, 4).Value = "auto"

# Row 4 (filled in after rows 3,5,6 - reuses C6orf174 string)
$ws.Cells.Item(4, 1).Value = "cg04580344"
$ws.Cells.Item(4, 2).Value = "C6orf174"
$ws.Cells.Item(4, 3).Value = "auto"
$ws.Cells.Item(4, 4).Value = "auto"

# Row 7
$ws.Cells.Item(7, 1).Value = "cg23928726"
$ws.Cells.Item(7, 2).Value = "PEX10"
$ws.Cells.Item(7, 3).Value = "auto"
$ws.Cells.Item(7, 4).Value = "auto"

$ws.Range("D7").Select()
